# Add a "Tap" column to the "Lines" sheet (between "Online" and "Notes"),
# and set the tap value of the transformer row ("Trafo1") to 0.
#
# This mirrors the commit: "'static' attributes are now accessible.
# Setting the tap turn via the input file is now also supported."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lines")

# --- Header row (row 1) ---------------------------------------------------
# Old layout: A..F = Name/From/To/Type/Length [km]/Online, G = Notes
# New layout: ... G = Tap (new), H = Notes (shifted right by one column)
$ws.Range("G1").Cut($ws.Range("H1"))
$ws.Range("F1").Copy($ws.Range("G1"))
$ws.Range("G1").Value = "Tap"

# --- "# Transformer" section header (row 2) --------------------------------
# G2 keeps its (blank, styled) look; duplicate the same styling into H2.
$ws.Range("G2").Copy($ws.Range("H2"))

# --- Trafo1 data row (row 3) -----------------------------------------------
# The old (empty, unstyled) G3 cell shifts to H3; the new G3 holds the tap
# value, styled like the other numeric cells in that row (e.g. F3).
$ws.Range("G3").Cut($ws.Range("H3"))
$ws.Range("F3").Copy($ws.Range("G3"))
$ws.Range("G3").Value = 0

# --- "# Lines" section header (row 4) --------------------------------------
# G4 keeps its (blank, styled) look; duplicate the same styling into H4.
$ws.Range("G4").Copy($ws.Range("H4"))

# --- Line data row (row 5) --------------------------------------------------
# Lines don't have a tap, so the old (empty, unstyled) G5 cell simply shifts
# to H5; no new value is written into G5.
$ws.Range("G5").Cut($ws.Range("H5"))

# Rows 6-8 (other line rows) never had a G/H cell before or after the edit,
# so they are intentionally left untouched.
